# Update "想去人数" (interested count) values in column F across sheets
# 展览 (rId1/sheet1), 本地生活 (rId3/sheet3) and 全部类型 (rId4/sheet4).
# 演出 (sheet2) is unchanged.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 21442
$ws1.Range("F3").Value  = 3393
$ws1.Range("F4").Value  = 860
$ws1.Range("F6").Value  = 560
$ws1.Range("F7").Value  = 816
$ws1.Range("F8").Value  = 306
$ws1.Range("F11").Value = 146
$ws1.Range("F12").Value = 595
$ws1.Range("F13").Value = 194
$ws1.Range("F15").Value = 45
$ws1.Range("F16").Value = 473
$ws1.Range("F17").Value = 241
$ws1.Range("F20").Value = 90
$ws1.Range("F21").Value = 165

# ---- Sheet: 本地生活 ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 736
$ws3.Range("F4").Value = 735
$ws3.Range("F5").Value = 1754

# ---- Sheet: 全部类型 ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 736
$ws4.Range("F4").Value  = 735
$ws4.Range("F5").Value  = 1754
$ws4.Range("F6").Value  = 21442
$ws4.Range("F7").Value  = 3393
$ws4.Range("F8").Value  = 860
$ws4.Range("F12").Value = 560
$ws4.Range("F13").Value = 817
$ws4.Range("F14").Value = 306
$ws4.Range("F20").Value = 146
$ws4.Range("F23").Value = 595
$ws4.Range("F25").Value = 194
$ws4.Range("F29").Value = 45
$ws4.Range("F30").Value = 473
$ws4.Range("F32").Value = 241
$ws4.Range("F37").Value = 90
$ws4.Range("F43").Value = 165
